$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (D value, E value). $null means "do not change".
$updates = @{
    2  = @("72.265.55", "  +1.14%  ")
    3  = @("2.669.30",  "  +1.61%  ")
    5  = @("599.06",    "  -1.37%  ")
    6  = @("176.49",    $null)
    7  = @($null,       "  +0.00%  ")
    8  = @($null,       "  -0.76%  ")
    9  = @("2.667.44",  "  +1.60%  ")
    10 = @("0.169",     "  +1.22%  ")
    11 = @($null,       "  +2.20%  ")
    12 = @($null,       "  +1.60%  ")
    13 = @("5.03",      "  +0.25%  ")
    14 = @("3.155.91",  $null)
    15 = @($null,       "  -1.70%  ")
    16 = @("72.073.70", "  +1.04%  ")
    17 = @("26.26",     "  -1.89%  ")
    18 = @("2.669.95",  "  +0.95%  ")
    19 = @("12.04",     "  +4.50%  ")
    20 = @("8.00",      "  +1.04%  ")
    21 = @("370.89",    "  -2.91%  ")
    22 = @("4.16",      "  +0.34%  ")
    23 = @($null,       "  +2.77%  ")
    24 = @("71.64",     "  -1.51%  ")
    25 = @($null,       "  +0.01%  ")
    26 = @("4.33",      "  -3.50%  ")
    27 = @("9.80",      "  +0.70%  ")
    28 = @("2.804.56",  "  +1.53%  ")
    29 = @("1.00",      "  -0.02%  ")
    30 = @("0.0₃0937",  "  -3.04%  ")
    31 = @("8.05",      "  -0.29%  ")
    32 = @("510.18",    "  -7.26%  ")
    33 = @($null,       "  -2.71%  ")
    34 = @($null,       "  -1.65%  ")
    35 = @($null,       "  +0.01%  ")
    36 = @("164.84",    "  -0.67%  ")
    37 = @("19.54",     "  +1.38%  ")
    38 = @($null,       "  +0.24%  ")
    39 = @($null,       "  -0.76%  ")
    40 = @($null,       "  -4.41%  ")
    41 = @("0.106",     "  -9.26%  ")
    42 = @($null,       "  +0.02%  ")
    43 = @("5.01",      "  -0.88%  ")
    44 = @($null,       "  -3.13%  ")
    45 = @("0.332",     "  -0.05%  ")
    46 = @("39.21",     "  -2.02%  ")
    47 = @("152.64",    "  -1.14%  ")
    48 = @($null,       "  +1.80%  ")
    49 = @($null,       "  +2.25%  ")
    50 = @($null,       "  +1.53%  ")
    51 = @("0.0768",    "  +1.80%  ")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $dVal = $pair[0]
    $eVal = $pair[1]

    if ($null -ne $dVal) {
        # Force the cell to remain plain text (many of these values look
        # like numbers, e.g. "599.06") while keeping the original,
        # un-styled General format -- matching how the source file stores
        # these as inline/shared strings without any numeric formatting.
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $dVal
        $cell.ClearFormats()
    }
    if ($null -ne $eVal) {
        $ws.Range("E$row").Value = $eVal
    }
}
